$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing row (user_id=0) ---
# place_id column (B) used to hold the email; it now says "admin"
$ws.Range("B2").Value = "admin"
# seat id changes to a new seat code
$ws.Range("C2").Value = "Main hall_0_main_3_2"
# refreshed timestamp
$ws.Range("D2").Value = 45500.23421097222

# --- Row 3: brand-new row (user_id=1) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "bogdan.yakupov@nu.edu.kz"
$ws.Range("C3").Value = "Main hall_0_main_3_4"
$ws.Range("D3").Value = 45500.23984773971
# copy row-2's formatting onto the new row so styles (border/bold/alignment on A, date format on D) match
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# --- Row 4: brand-new row (user_id=2) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "admin"
$ws.Range("C4").Value = "Main hall_0_main_0_3"
$ws.Range("D4").Value = 45500.24039318286
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$excel.CutCopyMode = 0
